$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 tweaks / deletions
$ws.Range("B2").Value = 435.435768
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 435.435768
$ws.Range("E2").Value = 516.92082400000004

# Row 3 tweaks / deletions
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 519.46723199999997
$ws.Range("D3").Value = 435.435768
$ws.Range("E3").Value = 519.46723199999997

# Update the selected range shown when the sheet is next opened
[void]$ws.Range("B1:E3").Select()
